$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7542783617973328
$ws.Range("B1").Value = 2.270795822143555
$ws.Range("D1").Value = 0.866680383682251
$ws.Range("E1").Value = 0.8072892427444458
